$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the stray note text in J25 ("Risultati davvero bassi con tempi di lavoro molto lunghi")
$ws.Range("J25").ClearContents()

# 2. Row 26 (MobileNetV2): bump Max epochs to 40 and fill in the measured results
$ws.Range("C26").Value = 40
$ws.Range("G26").Value = 0.96519999999999995
$ws.Range("H26").Value = 0.98080000000000001
$ws.Range("I26").Value = 0.97330000000000005

# 3. Add three new header+data(+blank) blocks at the bottom of the sheet, reusing the
#    existing formatted block (rows 28:30 = header / InceptionV3 data / blank spacer)
#    as the style template so the new cells come out with identical formatting.
$ws.Range("A28:I30").Copy($ws.Range("A40"))
$ws.Range("A28:I30").Copy($ws.Range("A43"))
$ws.Range("A28:I29").Copy($ws.Range("A46"))

# Block 1 (rows 40-42): GoogLeNet
$ws.Range("A41").Value = "GoogLeNet"
$ws.Range("B41").Value = "Adam"
$ws.Range("C41").Value = 80
$ws.Range("D41").Value = 0.0001
$ws.Range("E41").Value = 64
$ws.Range("F41").Value = 5
$ws.Range("G41").Value = 0.94779999999999998
$ws.Range("H41").Value = 0.95479999999999998
$ws.Range("I41").Value = 0.94330000000000003

# Block 2 (rows 43-45): ShuffleNet
$ws.Range("A44").Value = "ShuffleNet"
$ws.Range("B44").Value = "Adam"
$ws.Range("C44").Value = 100
$ws.Range("D44").Value = 0.0001
$ws.Range("E44").Value = 64
$ws.Range("F44").Value = 5
$ws.Range("G44").Value = 0.97389999999999999
$ws.Range("H44").Value = 0.98209999999999997
$ws.Range("I44").Value = 0.97499999999999998

# Block 3 (rows 46-47): SqueezeNet
$ws.Range("A47").Value = "SqueezeNet"
$ws.Range("B47").Value = "Adam"
$ws.Range("C47").Value = 100
$ws.Range("D47").Value = 0.0001
$ws.Range("E47").Value = 64
$ws.Range("F47").Value = 5
$ws.Range("G47").Value = 0.96519999999999995
$ws.Range("H47").Value = 0.9788
$ws.Range("I47").Value = 0.96989999999999998

# 4. Update the view: scroll so row 14 is at the top and select K40, matching the saved state
$ws.Range("K40").Select()
$ws.Application.ActiveWindow.ScrollRow = 14
